# Generate Report for Handoff
# The file f6c6938f-6784-4149-ac3a-94204185f2d8.md (row 3 on every sheet) has
# just been handed off for localization: its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", the handoff
# timestamps are refreshed, and the overview/per-locale sheets pick up a new
# "version mismatch" error detail message plus a wider Error Detail column.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/79f02a5898015e02b7ebba12531e7e3c1058d30d/e2e/f6c6938f-6784-4149-ac3a-94204185f2d8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d1162e18aaba57166ae9abb1a9a78443d31135f/e2e/f6c6938f-6784-4149-ac3a-94204185f2d8.md."

# --- Overview sheet: row 3 corresponds to f6c6938f-6784-4149-ac3a-94204185f2d8.md
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusReady
$overview.Range("F3").Value = $statusReady
$overview.Range("G3").Value = "2016-08-23 16:52:40"

# --- zh-cn sheet: row 3 is f6c6938f-6784-4149-ac3a-94204185f2d8.md
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusReady
$zhcn.Range("H3").Value = "2016-08-23 16:52:35"
$zhcn.Range("P3").Value = $errorMsg
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is f6c6938f-6784-4149-ac3a-94204185f2d8.md
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusReady
$dede.Range("H3").Value = "2016-08-23 16:52:40"
$dede.Range("P3").Value = $errorMsg
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
